$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '57.348.79'
$ws.Range('E2').Value = '  -0.51%  '

# Row 3
$ws.Range('D3').Value = '2.357.84'
$ws.Range('E3').Value = '  +1.10%  '

# Row 4
$ws.Range('E4').Value = '  +0.43%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.52'
$ws.Range('E5').Value = '  -0.26%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.58'
$ws.Range('E6').Value = '  +0.17%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.26%  '

# Row 8
$ws.Range('E8').Value = '  +0.35%  '

# Row 9
$ws.Range('E9').Value = '  -1.31%  '

# Row 10
$ws.Range('E10').Value = '  +5.17%  '

# Row 11
$ws.Range('E11').Value = '  -0.66%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('E12').Value = '  -0.81%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.38'
$ws.Range('E13').Value = '  +1.18%  '

# Row 14
$ws.Range('D14').Value = '2.781.08'
$ws.Range('E14').Value = '  +1.15%  '

# Row 15
$ws.Range('D15').Value = '57.359.32'
$ws.Range('E15').Value = '  +0.61%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000135'
$ws.Range('E16').Value = '  -0.58%  '

# Row 17
$ws.Range('D17').Value = '2.381.94'
$ws.Range('E17').Value = '  +1.71%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.61'
$ws.Range('E18').Value = '  -0.07%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '329.24'
$ws.Range('E19').Value = '  +1.92%  '

# Row 20
$ws.Range('E20').Value = '  -1.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  +1.00%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.11%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.30'
$ws.Range('E23').Value = '  +0.02%  '

# Row 24
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.72'
$ws.Range('E24').Value = '  +12.03%  '

# Row 25
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  +3.84%  '

# Row 26
$ws.Range('E26').Value = '  +0.36%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.34'
$ws.Range('E27').Value = '  +9.34%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0743'
$ws.Range('E28').Value = '  -0.91%  '

# Row 29
$ws.Range('E29').Value = '  -2.68%  '

# Row 30
$ws.Range('E30').Value = '  -0.17%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.28'
$ws.Range('E31').Value = '  -0.98%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.56'
$ws.Range('E32').Value = '  +0.71%  '

# Row 33
$ws.Range('E33').Value = '  +0.05%  '

# Row 34
$ws.Range('E34').Value = '  +2.17%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.07%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.928'
$ws.Range('E36').Value = '  -2.92%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.03'
$ws.Range('E37').Value = '  -0.66%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.62'
$ws.Range('E38').Value = '  +5.75%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.84'
$ws.Range('E39').Value = '  +3.32%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.84'
$ws.Range('E40').Value = '  +7.09%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.385'
$ws.Range('E41').Value = '  +0.65%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.66'
$ws.Range('E42').Value = '  +1.01%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.36'
$ws.Range('E43').Value = '  +3.09%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '283.96'
$ws.Range('E44').Value = '  +1.96%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0943'
$ws.Range('E45').Value = '  +1.30%  '

# Row 46
$ws.Range('E46').Value = '  -0.30%  '

# Row 47
$ws.Range('E47').Value = '  -0.24%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.22'
$ws.Range('E48').Value = '  +4.41%  '

# Row 49
$ws.Range('E49').Value = '  +1.25%  '

# Row 50
$ws.Range('B50').Value = 'Polygon'
$ws.Range('C50').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.387'
$ws.Range('E50').Value = '  +1.08%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.66'
$ws.Range('E51').Value = '  +3.33%  '
